$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and Column E (Volume(1h)) updates, keyed by row number.
$updates = @(
    @{ Row = 2;  D = "62.055.51";  E = "  +2.59%  " },
    @{ Row = 3;  D = "2.421.70";   E = "  +3.81%  " },
    @{ Row = 4;  E = "  +0.09%  " },
    @{ Row = 5;  D = "557.25";     E = "  +2.19%  " },
    @{ Row = 6;  D = "138.50";     E = "  +5.38%  " },
    @{ Row = 7;  E = "  +0.07%  " },
    @{ Row = 8;  D = "0.584";      E = "  +0.91%  " },
    @{ Row = 9;  D = "2.419.48";   E = "  +3.83%  " },
    @{ Row = 10; E = "  +2.57%  " },
    @{ Row = 11; D = "5.78";       E = "  +4.36%  " },
    @{ Row = 12; E = "  -0.02%  " },
    @{ Row = 13; E = "  +3.45%  " },
    @{ Row = 14; D = "25.76";      E = "  +8.79%  " },
    @{ Row = 15; D = "2.855.85";   E = "  +3.95%  " },
    @{ Row = 16; D = "62.018.04";  E = "  +2.59%  " },
    @{ Row = 17; D = "0.0000140";  E = "  +4.95%  " },
    @{ Row = 18; D = "2.424.23";   E = "  +3.91%  " },
    @{ Row = 19; D = "11.12";      E = "  +4.78%  " },
    @{ Row = 20; D = "343.80";     E = "  +9.05%  " },
    @{ Row = 21; E = "  +1.92%  " },
    @{ Row = 22; E = "  +2.54%  " },
    @{ Row = 23; E = "  +0.07%  " },
    @{ Row = 24; D = "65.04";      E = "  +1.62%  " },
    @{ Row = 25; E = "  -0.43%  " },
    @{ Row = 26; D = "1.00";       E = "  +0.25%  " },
    @{ Row = 27; B = "Fetch.AI"; C = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"; D = "1.53"; E = "  +12.26%  " },
    @{ Row = 28; B = "InternetComputer(DFINITY)"; C = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D = "8.28"; E = "  +5.27%  " },
    @{ Row = 29; D = "1.37";       E = "  +14.37%  " },
    @{ Row = 30; D = "0.0₃0786";   E = "  +7.09%  " },
    @{ Row = 31; D = "1.79";       E = "  +3.24%  " },
    @{ Row = 32; E = "  +6.25%  " },
    @{ Row = 33; D = "170.86" },
    @{ Row = 34; D = "1.44";       E = "  +4.69%  " },
    @{ Row = 35; D = "0.395";      E = "  +3.66%  " },
    @{ Row = 36; D = "373.83" },
    @{ Row = 37; D = "18.50";      E = "  +3.67%  " },
    @{ Row = 38; D = "4.49";       E = "  +10.01%  " },
    @{ Row = 40; E = "  -0.11%  " },
    @{ Row = 41; D = "1.68";       E = "  +9.32%  " },
    @{ Row = 42; D = "39.08";      E = "  +2.82%  " },
    @{ Row = 43; D = "146.02";     E = "  +5.92%  " },
    @{ Row = 44; D = "3.66";       E = "  +4.58%  " },
    @{ Row = 45; D = "20.75";      E = "  +7.88%  " },
    @{ Row = 47; D = "0.587";      E = "  +4.00%  " },
    @{ Row = 48; D = "0.0517";     E = "  +4.21%  " },
    @{ Row = 49; D = "18.03";      E = "  +6.22%  " },
    @{ Row = 50; E = "  +3.08%  " },
    @{ Row = 51; E = "  +3.28%  " }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey("B")) {
        $ws.Cells.Item($r, 2).Value = $u.B
    }
    if ($u.ContainsKey("C")) {
        $ws.Cells.Item($r, 3).Value = $u.C
    }
    if ($u.ContainsKey("D")) {
        # Price column holds text-formatted numbers (e.g. "62.055.51",
        # "557.25"). Force the cell to Text first so Excel's COM layer
        # doesn't silently coerce numeric-looking strings into real
        # numbers, then restore the default (unstyled) cell style.
        $cell = $ws.Cells.Item($r, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.Style = "Normal"
    }
    if ($u.ContainsKey("E")) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
